$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 05:44"

# Honduras (row 50)
$ws.Range("B50").Value = 67136
$ws.Range("C50").Value = 1087
$ws.Range("D50").Value = 17760
$ws.Range("E50").Value = 47311
$ws.Range("G50").Value = 7
$ws.Range("H50").Value = 2065

# Belice (row 160)
$ws.Range("B160").Value = 1458
$ws.Range("C160").Value = 23
$ws.Range("D160").Value = 458
$ws.Range("E160").Value = 981

# San Martin (Parte Holandesa) (row 173)
$ws.Range("B173").Value = 533
$ws.Range("C173").Value = 2
$ws.Range("D173").Value = 430
$ws.Range("E173").Value = 84

# Camboya (row 185)
$ws.Range("B185").Value = 275
$ws.Range("C185").Value = 1
$ws.Range("E185").Value = 1
